# Data Driven Testing added for User Sign Up Page....
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("I1").Value = "Zip code"
$ws.Range("J1").Value = "Mobile Number"

# New test-data emails (hyperlink formatting / relationships are preserved automatically)
$ws.Range("A2").Value = "tester987@gmail.com"
$ws.Range("A3").Value = "tester123@yahoomail.com"
$ws.Range("A4").Value = "tester1234@gmail.com"

# Whole data body switches to Text format so values (esp. the zip code /
# mobile number columns) are stored as literal strings, not numbers -- matches
# the new numFmtId="49" cellXfs entries.
$ws.Range("A2:K4").NumberFormat = "@"

$ws.Range("I2").Value = "12345"
$ws.Range("I3").Value = "98765"
$ws.Range("I4").Value = "54321"

$ws.Range("J2").Value = "6472152002"
$ws.Range("J3").Value = "6472152003"
$ws.Range("J4").Value = "6472152005"

# Move the active selection to A2 (was B15)
[void]$ws.Range("A2").Select()
